$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (G1:S1): insert "anybad" family of columns ---
$headerArr = New-Object "object[,]" 1,13
$headerArr[0,0] = "anybad"
$headerArr[0,1] = "rr_physical"
$headerArr[0,2] = "rr_stress"
$headerArr[0,3] = "rr_poverty"
$headerArr[0,4] = "rr_anybad"
$headerArr[0,5] = "rr_physical_95low"
$headerArr[0,6] = "rr_stress_95low"
$headerArr[0,7] = "rr_poverty_95low"
$headerArr[0,8] = "rr_anybad_95low"
$headerArr[0,9] = "rr_physical_95up"
$headerArr[0,10] = "rr_stress_95up"
$headerArr[0,11] = "rr_poverty_95up"
$headerArr[0,12] = "rr_anybad_95up"
$ws.Range("G1:S1").Value = $headerArr

# --- Data rows (G2:S19) ---
$dataArr = New-Object "object[,]" 18,13
$dataArr[0,0] = 0.8175230068533549
$dataArr[0,1] = 1.044098255695232
$dataArr[0,2] = 1.237123360808975
$dataArr[0,3] = 1.050483634589852
$dataArr[0,4] = 1.053956506267101
$dataArr[0,5] = 1.013214108030252
$dataArr[0,6] = 1.177261567584297
$dataArr[0,7] = 0.9393582234967901
$dataArr[0,8] = 1.032635334999701
$dataArr[0,9] = 1.080641699909626
$dataArr[0,10] = 1.308272987123011
$dataArr[0,11] = 1.177921432600948
$dataArr[0,12] = 1.080691074729922
$dataArr[1,0] = 0.8386909240581151
$dataArr[1,1] = 1.054004392494821
$dataArr[1,2] = 1.193829524087098
$dataArr[1,3] = 1.106747649520144
$dataArr[1,4] = 1.05346833756508
$dataArr[1,5] = 1.021972262622803
$dataArr[1,6] = 1.146143366440801
$dataArr[1,7] = 1.0015820107145
$dataArr[1,8] = 1.035052565250807
$dataArr[1,9] = 1.090851643787171
$dataArr[1,10] = 1.250298774655401
$dataArr[1,11] = 1.214972318295113
$dataArr[1,12] = 1.077299313986833
$dataArr[2,0] = 0.8193264327177946
$dataArr[2,1] = 1.058382813377533
$dataArr[2,2] = 1.211358862466423
$dataArr[2,3] = 1.045574404598978
$dataArr[2,4] = 1.05342325241806
$dataArr[2,5] = 1.016730250379123
$dataArr[2,6] = 1.158456093903216
$dataArr[2,7] = 0.9450488182099978
$dataArr[2,8] = 1.032529503223918
$dataArr[2,9] = 1.103426121425787
$dataArr[2,10] = 1.272277271588753
$dataArr[2,11] = 1.158013837665733
$dataArr[2,12] = 1.076964136236134
$dataArr[3,0] = 0.8369997955036165
$dataArr[3,1] = 1.069159225410057
$dataArr[3,2] = 1.19401426989448
$dataArr[3,3] = 1.098817005782724
$dataArr[3,4] = 1.054028887750431
$dataArr[3,5] = 1.027933323073099
$dataArr[3,6] = 1.146196841934884
$dataArr[3,7] = 1.001467766537552
$dataArr[3,8] = 1.036108510667731
$dataArr[3,9] = 1.112353739243522
$dataArr[3,10] = 1.249971289715862
$dataArr[3,11] = 1.201077429356429
$dataArr[3,12] = 1.074023842891252
$dataArr[4,0] = 0.8268440813116509
$dataArr[4,1] = 1.060409947850499
$dataArr[4,2] = 1.161259776203038
$dataArr[4,3] = 1.056943099933577
$dataArr[4,4] = 1.051200363677547
$dataArr[4,5] = 1.017066322711784
$dataArr[4,6] = 1.117767919908355
$dataArr[4,7] = 0.9305699109706342
$dataArr[4,8] = 1.031149756242346
$dataArr[4,9] = 1.107693991981422
$dataArr[4,10] = 1.207719926536266
$dataArr[4,11] = 1.207015715922681
$dataArr[4,12] = 1.075558763509092
$dataArr[5,0] = 0.8398475075689038
$dataArr[5,1] = 1.072862908953344
$dataArr[5,2] = 1.137085132776632
$dataArr[5,3] = 1.122686879311817
$dataArr[5,4] = 1.053084970434521
$dataArr[5,5] = 1.029044207173212
$dataArr[5,6] = 1.102392628159572
$dataArr[5,7] = 1.001829012983137
$dataArr[5,8] = 1.035097045268287
$dataArr[5,9] = 1.120741644427255
$dataArr[5,10] = 1.177113732523665
$dataArr[5,11] = 1.25171213644987
$dataArr[5,12] = 1.075633723840007
$dataArr[6,0] = 0.7686456608164728
$dataArr[6,1] = 1.064504197049108
$dataArr[6,2] = 1.201974245403202
$dataArr[6,3] = 1.056139266225322
$dataArr[6,4] = 1.068409017689398
$dataArr[6,5] = 1.019350671261173
$dataArr[6,6] = 1.150516758427119
$dataArr[6,7] = 0.9320771727041821
$dataArr[6,8] = 1.043327013531331
$dataArr[6,9] = 1.110405551463595
$dataArr[6,10] = 1.260815716801262
$dataArr[6,11] = 1.20311753225668
$dataArr[6,12] = 1.09744546864995
$dataArr[7,0] = 0.7679803848535932
$dataArr[7,1] = 1.086727703560965
$dataArr[7,2] = 1.169755052274979
$dataArr[7,3] = 1.120523313161906
$dataArr[7,4] = 1.076906417273369
$dataArr[7,5] = 1.035001138198733
$dataArr[7,6] = 1.129628157823973
$dataArr[7,7] = 1.001784824149812
$dataArr[7,8] = 1.051708400033905
$dataArr[7,9] = 1.143951573749473
$dataArr[7,10] = 1.218419535632835
$dataArr[7,11] = 1.245031906178167
$dataArr[7,12] = 1.103342044859157
$dataArr[8,0] = 0.7784494204094269
$dataArr[8,1] = 1.074349537467539
$dataArr[8,2] = 1.177110981494486
$dataArr[8,3] = 1.05442242154026
$dataArr[8,4] = 1.065510150238786
$dataArr[8,5] = 1.021298947293009
$dataArr[8,6] = 1.132783373620724
$dataArr[8,7] = 0.9341703872701977
$dataArr[8,8] = 1.042270096761871
$dataArr[8,9] = 1.130754085254148
$dataArr[8,10] = 1.229568451598431
$dataArr[8,11] = 1.196362850826908
$dataArr[8,12] = 1.092956980082138
$dataArr[9,0] = 0.7735265738580027
$dataArr[9,1] = 1.096142850300793
$dataArr[9,2] = 1.170936579646361
$dataArr[9,3] = 1.118150736159423
$dataArr[9,4] = 1.075068048885502
$dataArr[9,5] = 1.038417612170542
$dataArr[9,6] = 1.128334827182179
$dataArr[9,7] = 1.001744630190954
$dataArr[9,8] = 1.051408116696234
$dataArr[9,9] = 1.156090774583735
$dataArr[9,10] = 1.223374416771505
$dataArr[9,11] = 1.240218270360443
$dataArr[9,12] = 1.10412702004948
$dataArr[10,0] = 0.7782213785350327
$dataArr[10,1] = 1.082125677360067
$dataArr[10,2] = 1.148395490054291
$dataArr[10,3] = 1.05956248074361
$dataArr[10,4] = 1.065577579795864
$dataArr[10,5] = 1.023281762176376
$dataArr[10,6] = 1.11136246562521
$dataArr[10,7] = 0.9279531562222474
$dataArr[10,8] = 1.041697088410233
$dataArr[10,9] = 1.143563474572525
$dataArr[10,10] = 1.192031752212329
$dataArr[10,11] = 1.216672860596275
$dataArr[10,12] = 1.092734432266782
$dataArr[11,0] = 0.7678835954183992
$dataArr[11,1] = 1.105665731865027
$dataArr[11,2] = 1.133812776795585
$dataArr[11,3] = 1.129839062148083
$dataArr[11,4] = 1.076938499598331
$dataArr[11,5] = 1.04320478350242
$dataArr[11,6] = 1.102747176908807
$dataArr[11,7] = 1.001918740400947
$dataArr[11,8] = 1.053282848891016
$dataArr[11,9] = 1.173172326714252
$dataArr[11,10] = 1.170634527022718
$dataArr[11,11] = 1.263791016657192
$dataArr[11,12] = 1.103020414838134
$dataArr[12,0] = 0.705501664307754
$dataArr[12,1] = 1.110297834755405
$dataArr[12,2] = 1.155644161590905
$dataArr[12,3] = 1.059093269015841
$dataArr[12,4] = 1.08708002593324
$dataArr[12,5] = 1.032028395153306
$dataArr[12,6] = 1.116487137873932
$dataArr[12,7] = 0.9286757853284469
$dataArr[12,8] = 1.055772744729348
$dataArr[12,9] = 1.194189475733554
$dataArr[12,10] = 1.203620956084896
$dataArr[12,11] = 1.214669099716139
$dataArr[12,12] = 1.124642813695595
$dataArr[13,0] = 0.7392574504319697
$dataArr[13,1] = 1.120484927432569
$dataArr[13,2] = 1.12390602488182
$dataArr[13,3] = 1.12994412284401
$dataArr[13,4] = 1.086427069130977
$dataArr[13,5] = 1.049752996338377
$dataArr[13,6] = 1.092596943074343
$dataArr[13,7] = 1.001921930265669
$dataArr[13,8] = 1.058579656712467
$dataArr[13,9] = 1.199021002126987
$dataArr[13,10] = 1.161494096996591
$dataArr[13,11] = 1.262459936810329
$dataArr[13,12] = 1.11816940764967
$dataArr[14,0] = 0.7491129596126874
$dataArr[14,1] = 1.099664863840935
$dataArr[14,2] = 1.148532833266064
$dataArr[14,3] = 1.05639414197546
$dataArr[14,4] = 1.074184629708982
$dataArr[14,5] = 1.027238502104199
$dataArr[14,6] = 1.109570089706756
$dataArr[14,7] = 0.9323585693685164
$dataArr[14,8] = 1.045488445675182
$dataArr[14,9] = 1.176537087374862
$dataArr[14,10] = 1.197868008031604
$dataArr[14,11] = 1.199512780459306
$dataArr[14,12] = 1.110147578992588
$dataArr[15,0] = 0.7751774536136742
$dataArr[15,1] = 1.106174303322158
$dataArr[15,2] = 1.139799305085285
$dataArr[15,3] = 1.125287773663258
$dataArr[15,4] = 1.074520839774421
$dataArr[15,5] = 1.044191218642741
$dataArr[15,6] = 1.10149487311788
$dataArr[15,7] = 1.001861028740489
$dataArr[15,8] = 1.048967267030918
$dataArr[15,9] = 1.175625242699108
$dataArr[15,10] = 1.186480890678051
$dataArr[15,11] = 1.255158002195386
$dataArr[15,12] = 1.110336569341654
$dataArr[16,0] = 0.7507024640566946
$dataArr[16,1] = 1.118240713509426
$dataArr[16,2] = 1.117002258861124
$dataArr[16,3] = 1.060457427763088
$dataArr[16,4] = 1.073714630148952
$dataArr[16,5] = 1.033145783033963
$dataArr[16,6] = 1.088587579263204
$dataArr[16,7] = 0.9269355032269437
$dataArr[16,8] = 1.047706556614987
$dataArr[16,9] = 1.207720122488238
$dataArr[16,10] = 1.151406130392854
$dataArr[16,11] = 1.219519675521085
$dataArr[16,12] = 1.103358078668907
$dataArr[17,0] = 0.7712961313750915
$dataArr[17,1] = 1.130363315869883
$dataArr[17,2] = 1.101589580515479
$dataArr[17,3] = 1.132914905882109
$dataArr[17,4] = 1.075807362844742
$dataArr[17,5] = 1.052815520546763
$dataArr[17,6] = 1.077949976174376
$dataArr[17,7] = 1.001960602022187
$dataArr[17,8] = 1.051911592115921
$dataArr[17,9] = 1.216269105300662
$dataArr[17,10] = 1.128794505774226
$dataArr[17,11] = 1.270580219275513
$dataArr[17,12] = 1.102033319354732
$ws.Range("G2:S19").Value = $dataArr

# Header style: bold+center, matching other header cells
$ws.Range("G1:S1").Font.Bold = $true
$ws.Range("G1:S1").HorizontalAlignment = -4108  # xlCenter